$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(206).Insert()

$ws.Range("A206").Value = 11
$ws.Range("B206").Value = "Vega Monumental Concepción"
$ws.Range("C206").Value = "Bíobío"
$ws.Range("D206").Value = 44918
$ws.Range("D206").NumberFormat = $ws.Range("D207").NumberFormat
$ws.Range("E206").Value = 8
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100108
$ws.Range("H206").Value = "Tropicales y subtropicales"
$ws.Range("I206").Value = 100108005
$ws.Range("J206").Value = "Piña"
$ws.Range("K206").Value = "Caramelo"
$ws.Range("L206").Value = "Segunda"
$ws.Range("M206").Value = 100
$ws.Range("N206").Value = 21000
$ws.Range("O206").Value = 22000
$ws.Range("P206").Value = 21500
$ws.Range("Q206").Value = "$/caja 14 unidades"
$ws.Range("R206").Value = "Ecuador"
$ws.Range("S206").Value = 1536
$ws.Range("T206").Value = 14
